# Update "paises" (countries) sheet: refresh shared-string order for the
# reordered country-name rows and write the latest case counts, matching
# the 18:20 data refresh from the 17:50 snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Country names that moved to a new row position ---
$ws.Range("A1").Value = "Datos actualizados a 29 de Marzo de 2020 a las 18:20"
$ws.Range("A63").Value = "Argelia"
$ws.Range("A64").Value = "Barein"
$ws.Range("A65").Value = "Emiratos Arabes Unidos"
$ws.Range("A81").Value = "Moldavia"
$ws.Range("A82").Value = "Republica de Macedonia"
$ws.Range("A83").Value = "Kuwait"
$ws.Range("A84").Value = "Kazajistan"
$ws.Range("A85").Value = "Jordania"
$ws.Range("A146").Value = "Niger"
$ws.Range("A147").Value = "Mali"
$ws.Range("A155").Value = "Mongolia"
$ws.Range("A156").Value = "Guinea Ecuatorial"
$ws.Range("A159").Value = "Bahamas"
$ws.Range("A160").Value = "Namibia"
$ws.Range("A164").Value = "Siria"
$ws.Range("A165").Value = "Laos"
$ws.Range("A166").Value = "Seychelles"
$ws.Range("A167").Value = "Birmania"
$ws.Range("A168").Value = "Surinam"
$ws.Range("A169").Value = "Mozambique"
$ws.Range("A170").Value = "Libia"
$ws.Range("A171").Value = "Guyana"
$ws.Range("A172").Value = "Islas Caimanes"
$ws.Range("A173").Value = "Curazao"
$ws.Range("A174").Value = "Antigua y Barbuda"
$ws.Range("A176").Value = "Gabon"
$ws.Range("A177").Value = "Benin"
$ws.Range("A178").Value = "Santa Sede"
$ws.Range("A180").Value = "San Martin (Parte Holandesa)"
$ws.Range("A183").Value = "Fiyi"
$ws.Range("A184").Value = "San Bartolome"
$ws.Range("A185").Value = "Angola"
$ws.Range("A186").Value = "Mauritania"
$ws.Range("A187").Value = "Sudan"
$ws.Range("A188").Value = "Nepal"
$ws.Range("A189").Value = "Islas Turcas y Caicos"
$ws.Range("A190").Value = "Butan"
$ws.Range("A191").Value = "Santa Lucia"
$ws.Range("A192").Value = "Nicaragua"
$ws.Range("A193").Value = "Republica del Chad"
$ws.Range("A196").Value = "Somalia"

# --- Updated numeric data (latest case counts) ---
$ws.Range("B4").Value = 131366
$ws.Range("C4").Value = 7788
$ws.Range("D4").Value = 4378
$ws.Range("E4").Value = 124660
$ws.Range("F4").Value = 2948
$ws.Range("G4").Value = 107
$ws.Range("H4").Value = 2328
$ws.Range("B5").Value = 97689
$ws.Range("C5").Value = 5217
$ws.Range("D5").Value = 13030
$ws.Range("E5").Value = 73880
$ws.Range("F5").Value = 3906
$ws.Range("G5").Value = 756
$ws.Range("H5").Value = 10779
$ws.Range("D13").Value = 250
$ws.Range("E13").Value = 9845
$ws.Range("F13").Value = 972
$ws.Range("E26").Value = 2047
$ws.Range("G26").Value = 8
$ws.Range("H26").Value = 35
$ws.Range("D29").Value = 75
$ws.Range("E29").Value = 2057
$ws.Range("G29").Value = 1
$ws.Range("H29").Value = 7
$ws.Range("B63").Value = 511
$ws.Range("C63").Value = 57
$ws.Range("D63").Value = 31
$ws.Range("E63").Value = 451
$ws.Range("F63").Value = 0
$ws.Range("H63").Value = 29
$ws.Range("B64").Value = 499
$ws.Range("C64").Value = 23
$ws.Range("D64").Value = 272
$ws.Range("E64").Value = 223
$ws.Range("F64").Value = 1
$ws.Range("H64").Value = 4
$ws.Range("B65").Value = 468
$ws.Range("D65").Value = 55
$ws.Range("E65").Value = 411
$ws.Range("F65").Value = 2
$ws.Range("H65").Value = 2
$ws.Range("B81").Value = 263
$ws.Range("C81").Value = 32
$ws.Range("D81").Value = 2
$ws.Range("E81").Value = 259
$ws.Range("F81").Value = 33
$ws.Range("G81").Value = 0
$ws.Range("H81").Value = 2
$ws.Range("B82").Value = 259
$ws.Range("C82").Value = 18
$ws.Range("D82").Value = 3
$ws.Range("E82").Value = 250
$ws.Range("F82").Value = 1
$ws.Range("G82").Value = 2
$ws.Range("H82").Value = 6
$ws.Range("B83").Value = 255
$ws.Range("C83").Value = 20
$ws.Range("D83").Value = 67
$ws.Range("E83").Value = 188
$ws.Range("F83").Value = 12
$ws.Range("H83").Value = 0
$ws.Range("B84").Value = 251
$ws.Range("C84").Value = 23
$ws.Range("E84").Value = 232
$ws.Range("F84").Value = 0
$ws.Range("B85").Value = 246
$ws.Range("D85").Value = 18
$ws.Range("F85").Value = 3
$ws.Range("H85").Value = 1
$ws.Range("C146").Value = 8
$ws.Range("C147").Value = 0
$ws.Range("C159").Value = 1
$ws.Range("D159").Value = 1
$ws.Range("E159").Value = 10
$ws.Range("B160").Value = 11
$ws.Range("C160").Value = 3
$ws.Range("D160").Value = 2
$ws.Range("B164").Value = 9
$ws.Range("C164").Value = 4
$ws.Range("G164").Value = 1
$ws.Range("H164").Value = 1
$ws.Range("C169").Value = 0
$ws.Range("C170").Value = 5
$ws.Range("E170").Value = 8
$ws.Range("H170").Value = 0
$ws.Range("D172").Value = 0
$ws.Range("E172").Value = 7
$ws.Range("B173").Value = 8
$ws.Range("D173").Value = 2
$ws.Range("E173").Value = 5
$ws.Range("H173").Value = 1
$ws.Range("E174").Value = 7
$ws.Range("H174").Value = 0
$ws.Range("B176").Value = 7
$ws.Range("H176").Value = 1
$ws.Range("C180").Value = 0
$ws.Range("E180").Value = 6
$ws.Range("G180").Value = 0
$ws.Range("H180").Value = 0
$ws.Range("D187").Value = 0
$ws.Range("H187").Value = 1
$ws.Range("D188").Value = 1
$ws.Range("H188").Value = 0
$ws.Range("C189").Value = 0
$ws.Range("C190").Value = 1
$ws.Range("C191").Value = 1
$ws.Range("D191").Value = 1
$ws.Range("H191").Value = 0
$ws.Range("C192").Value = 0
$ws.Range("D192").Value = 0
